$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# Columns D..AJ map to spreadsheet columns 4..36.
# Row 2 = 2014/12, Row 3 = 2015/12, Row 4 = 2016/12, Row 5 = 2017/12, Row 6 = 2018/12.
# For each of those rows the financial figures were corrected (previous values
# were off by orders of magnitude / wrong units). Rows 7-9 (the forecast years
# 2019/12(E), 2020/12(E), 2021/12(E)) had their estimate data removed entirely,
# leaving only the index/label columns A-C.

function Set-RowValues($row, $startCol, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $val = $values[$i]
        if ($val -ne $null) {
            $ws.Cells.Item($row, $startCol + $i).Value = $val
        }
    }
}

# Row 2 (2014/12) - columns D(4) through AJ(36)
Set-RowValues 2 4 @(8026,82,82,-22,-23,-22,-1,4386,2704,1682,1663,19,300,201,-191,-47,194,7,1465,1.02,-0.29,-1.29,-0.52,160.77,453.92,-362,-37,27722,0.48,0,0,0,6000000)

# Row 3 (2015/12) - columns D(4) through AJ(36)
Set-RowValues 3 4 @(8016,231,231,121,78,94,-16,4319,2562,1757,1753,3,300,281,-69,-25,70,211,1464,2.89,0.98,5.53,1.8,145.85,483.32,1573,11.38,29222,0.61,750,4.19,47.67,6000000)

# Row 4 (2016/12) - columns D(4) through AJ(36)
Set-RowValues 4 4 @(7980,411,411,391,314,311,3,4058,2008,2050,2043,6,300,407,193,-655,82,325,861,5.15,3.94,16.39,7.5,97.97,579.78,5186,5.73,34051,0.87,1000,3.37,19.28,6000000)

# Row 5 (2017/12) - columns D(4) through AJ(36)
Set-RowValues 5 4 @(9121,244,244,223,175,180,-5,4274,2134,2141,2139,1,300,270,-284,2,264,6,904,2.67,1.92,8.6,4.21,99.68000000000001,612.23,2999,9.5,35655,0.8,750,2.63,25.01,6000000)

# Row 6 (2018/12) - columns D(4) through AJ(36); J6 and O6 have no source
# data (they were blank before the edit too) so they are skipped via $null.
Set-RowValues 6 4 @(9403,178,178,142,108,113,$null,4542,2355,2187,2191,$null,300,143,-386,149,388,-245,1104,1.89,1.15,5.24,2.45,107.69,629.38,1891,10.94,36515,0.57,500,2.42,26.44,6000000)

# Rows 7, 8, 9 (forecast years 2019/12(E), 2020/12(E), 2021/12(E)):
# all financial data columns (D..AJ) are removed, only A-C remain.
$ws.Range("D7:AJ9").ClearContents()
